$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 846.3692
$ws.Range("J17").Value = 849.40985
$ws.Range("L17").Value = 2548.22955
$ws.Range("N17").Value = -2884.22955
$ws.Range("H18").Value = 940
$ws.Range("I18").Value = 940
$ws.Range("K18").Value = 940
$ws.Range("M18").Value = -656
$ws.Range("H19").Value = 2426.2144
$ws.Range("I19").Value = 2006.6
$ws.Range("J19").Value = 3475.25
$ws.Range("K19").Value = 2006.6
$ws.Range("L19").Value = 3475.25
$ws.Range("M19").Value = -1831.6
$ws.Range("N19").Value = -3825.25
$ws.Range("H98").Value = 3552.2222
$ws.Range("I98").Value = 3367.2856
$ws.Range("K98").Value = 3367.2856
$ws.Range("M98").Value = -1869.2856
$ws.Range("H107").Value = 435
$ws.Range("I107").Value = 458.53845
$ws.Range("K107").Value = 458.53845
$ws.Range("M107").Value = 1461.46155
$ws.Range("H114").Value = 68999.5
$ws.Range("I114").Value = 68999.5
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 68999.5
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = -64660.5
$ws.Range("N114").ClearContents()
$ws.Range("H122").Value = 3552.2222
$ws.Range("I122").Value = 3367.2856
$ws.Range("K122").Value = 10101.8568
$ws.Range("M122").Value = -7651.856800000001
$ws.Range("H129").Value = 1222.875
$ws.Range("I129").Value = 1018.75
$ws.Range("J129").Value = 1427
$ws.Range("K129").Value = 3056.25
$ws.Range("L129").Value = 4281
$ws.Range("M129").Value = 1943.75
$ws.Range("N129").Value = -14281
$ws.Range("H132").Value = 2550.6
$ws.Range("I132").Value = 2550.6
$ws.Range("K132").Value = 7651.799999999999
$ws.Range("M132").Value = -5121.799999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1617.2572
$ws.Range("I32").Value = 1617.2572
$ws.Range("K32").Value = 1617.2572
$ws.Range("M32").Value = -1330.2572
$ws.Range("H61").Value = 4668.4546
$ws.Range("I61").Value = 4194.143
$ws.Range("K61").Value = 4194.143
$ws.Range("M61").Value = -3982.143
$ws.Range("H94").Value = 34971.75
$ws.Range("J94").Value = 34971.75
$ws.Range("L94").Value = 34971.75
$ws.Range("N94").Value = -36773.75
$ws.Range("H110").Value = 2349.862
$ws.Range("I110").Value = 1831.9
$ws.Range("J110").Value = 3500.889
$ws.Range("K110").Value = 1831.9
$ws.Range("L110").Value = 3500.889
$ws.Range("M110").Value = 213.0999999999999
$ws.Range("N110").Value = -7590.889
$ws.Range("H132").Value = 15630274
$ws.Range("I132").Value = 4845.0835
$ws.Range("J132").Value = 62506560
$ws.Range("K132").Value = 14535.2505
$ws.Range("L132").Value = 187519680
$ws.Range("M132").Value = -12005.2505
$ws.Range("N132").Value = -187524740
$ws.Range("H136").Value = 4668.4546
$ws.Range("I136").Value = 4194.143
$ws.Range("K136").Value = 12582.429
$ws.Range("M136").Value = -10032.429

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 203.81818
$ws.Range("J80").Value = 241.57143
$ws.Range("L80").Value = 241.57143
$ws.Range("N80").Value = -2237.57143
$ws.Range("H83").Value = 203.81818
$ws.Range("J83").Value = 241.57143
$ws.Range("L83").Value = 1207.85715
$ws.Range("N83").Value = -11191.85715
$ws.Range("H86").Value = 4005.1667
$ws.Range("I86").Value = 2682
$ws.Range("K86").Value = 2682
$ws.Range("M86").Value = -1559
$ws.Range("H89").Value = 4005.1667
$ws.Range("I89").Value = 2682
$ws.Range("K89").Value = 13410
$ws.Range("M89").Value = -7794
$ws.Range("H99").Value = 3864.2273
$ws.Range("I99").Value = 3868.0625
$ws.Range("J99").Value = 3854
$ws.Range("K99").Value = 3868.0625
$ws.Range("L99").Value = 3854
$ws.Range("M99").Value = -2370.0625
$ws.Range("N99").Value = -6850
$ws.Range("H134").Value = 23818172
$ws.Range("I134").Value = 11137.4
$ws.Range("J134").Value = 83335760
$ws.Range("K134").Value = 33412.2
$ws.Range("L134").Value = 250007280
$ws.Range("M134").Value = -30877.2
$ws.Range("N134").Value = -250012350

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2240.0967
$ws.Range("I31").Value = 1386.9
$ws.Range("J31").Value = 3791.3635
$ws.Range("K31").Value = 1386.9
$ws.Range("L31").Value = 3791.3635
$ws.Range("M31").Value = -1091.9
$ws.Range("N31").Value = -4381.363499999999
$ws.Range("H34").Value = 2240.0967
$ws.Range("I34").Value = 1386.9
$ws.Range("J34").Value = 3791.3635
$ws.Range("K34").Value = 1386.9
$ws.Range("L34").Value = 3791.3635
$ws.Range("M34").Value = -1184.9
$ws.Range("N34").Value = -4195.363499999999
$ws.Range("H58").Value = 2252.375
$ws.Range("I58").Value = 2000
$ws.Range("J58").Value = 2288.4285
$ws.Range("K58").Value = 2000
$ws.Range("L58").Value = 2288.4285
$ws.Range("M58").Value = -1797
$ws.Range("N58").Value = -2694.4285
$ws.Range("H99").Value = 4906.9165
$ws.Range("I99").Value = 4877.4443
$ws.Range("J99").Value = 4995.3335
$ws.Range("K99").Value = 4877.4443
$ws.Range("L99").Value = 4995.3335
$ws.Range("M99").Value = -3379.4443
$ws.Range("N99").Value = -7991.3335
$ws.Range("H126").Value = 4906.9165
$ws.Range("I126").Value = 4877.4443
$ws.Range("J126").Value = 4995.3335
$ws.Range("K126").Value = 14632.3329
$ws.Range("L126").Value = 14986.0005
$ws.Range("M126").Value = -12162.3329
$ws.Range("N126").Value = -19926.0005
$ws.Range("H132").Value = 5632.75
$ws.Range("I132").Value = 5453.2856
$ws.Range("K132").Value = 16359.8568
$ws.Range("M132").Value = -13829.8568
$ws.Range("H133").Value = 49155
$ws.Range("J133").Value = 49155
$ws.Range("L133").Value = 49155
$ws.Range("N133").Value = -54215
$ws.Range("H134").Value = 8337068.5
$ws.Range("I134").Value = 2981.4443
$ws.Range("K134").Value = 8944.332900000001
$ws.Range("M134").Value = -6409.332900000001
$ws.Range("H136").Value = 2252.375
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 2288.4285
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 6865.2855
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -11965.2855

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1680.4
$ws.Range("I11").Value = 600
$ws.Range("K11").Value = 1800
$ws.Range("M11").Value = -1660
$ws.Range("H113").Value = 657.5454999999999
$ws.Range("J113").Value = 840
$ws.Range("L113").Value = 2520
$ws.Range("N113").Value = -6860
$ws.Range("H132").Value = 1274.75
$ws.Range("I132").Value = 700
$ws.Range("J132").Value = 1466.3334
$ws.Range("K132").Value = 6300
$ws.Range("L132").Value = 13197.0006
$ws.Range("N132").Value = -18257.0006
$ws.Range("M132").Value = -3770

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3135.4
$ws.Range("I102").Value = 2448
$ws.Range("J102").Value = 3593.6667
$ws.Range("K102").Value = 2448
$ws.Range("L102").Value = 3593.6667
$ws.Range("M102").Value = -826
$ws.Range("N102").Value = -6837.6667
$ws.Range("H132").Value = 3222
$ws.Range("I132").Value = 3166.5
$ws.Range("J132").Value = 3333
$ws.Range("K132").Value = 9499.5
$ws.Range("L132").Value = 9999
$ws.Range("M132").Value = -6969.5
$ws.Range("N132").Value = -15059
$ws.Range("H135").Value = 48712.5
$ws.Range("I135").Value = 49778
$ws.Range("J135").Value = 48357.332
$ws.Range("K135").Value = 49778
$ws.Range("L135").Value = 48357.332
$ws.Range("N135").Value = -58497.332
$ws.Range("M135").Value = -44708

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3585.9
$ws.Range("I61").Value = 3076.25
$ws.Range("K61").Value = 3076.25
$ws.Range("M61").Value = -2874.25
$ws.Range("H113").Value = 3585.9
$ws.Range("I113").Value = 3076.25
$ws.Range("K113").Value = 3076.25
$ws.Range("M113").Value = -906.25
$ws.Range("H122").Value = 3439.261
$ws.Range("I122").Value = 2886.3333
$ws.Range("J122").Value = 3522.2
$ws.Range("K122").Value = 8658.999899999999
$ws.Range("L122").Value = 10566.6
$ws.Range("M122").Value = -6208.999899999999
$ws.Range("N122").Value = -15466.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 18000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 18000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 18000
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -19262
$ws.Range("H61").Value = 13663
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 13663
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 13663
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -14247
$ws.Range("H110").Value = 60625
$ws.Range("J110").Value = 60625
$ws.Range("L110").Value = 60625
$ws.Range("N110").Value = -68805

